$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.961716651916504
$ws.Range("B1").Value = 2.659371376037598
$ws.Range("C1").Value = 1.905655741691589
$ws.Range("D1").Value = 1.742602229118347
$ws.Range("E1").Value = 1.751873254776001
